$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44628
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100104
$ws.Cells.Item($row, 8).Value = "Frutos de pepita"
$ws.Cells.Item($row, 9).Value = 100104003
$ws.Cells.Item($row, 10).Value = "Membrillo"
$ws.Cells.Item($row, 11).Value = "Champion"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 14
$ws.Cells.Item($row, 14).Value = 400000
$ws.Cells.Item($row, 15).Value = 410000
$ws.Cells.Item($row, 16).Value = 405000
$ws.Cells.Item($row, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 900
$ws.Cells.Item($row, 20).Value = 450
